# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet, insert a new (blank) column at N so the
# existing "Late" / "Outstanding" columns (previously N & P, with a spacer
# blank column O) shift one place to the right (to O & Q, with the new
# spacer blank column now at N). Give the freshly inserted column the same
# display width as it had in the saved workbook, then make "Repayment
# Schedule" the active sheet/tab with S9 selected (replacing "Transactions"
# as the active tab).

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N (shifts N:P -> O:Q).
$schedule.Columns("N:N").Insert()

# Match the column width recorded for the newly inserted column.
$schedule.Columns("N:N").ColumnWidth = 9.166666666666666

# Make "Repayment Schedule" the active sheet and select cell S9 on it
# (this also clears "Transactions" as the active/selected tab).
$schedule.Activate()
$schedule.Range("S9").Select() | Out-Null
